$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.002.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.262.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.81%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.646"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.89%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.70%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.450"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.89%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.72%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.35%  "

# Row 13
$ws.Range("E13").Value = "  +1.92%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.597.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.92%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.83%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.836"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.99%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.261.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.921.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.96%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "

# Row 22
$ws.Range("E22").Value = "  -1.38%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("E25").Value = "  -4.18%  "

# Row 26
$ws.Range("E26").Value = "  +21.74%  "

# Row 27
$ws.Range("E27").Value = "  -4.71%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.87%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.65%  "

# Row 31
$ws.Range("E31").Value = "  +0.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.48%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.20%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.35%  "

# Row 35
$ws.Range("E35").Value = "  -1.60%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.44%  "

# Row 37
$ws.Range("E37").Value = "  -2.56%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.00%  "

# Row 39
$ws.Range("E39").Value = "  -2.81%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0256"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.15%  "

# Row 41
$ws.Range("E41").Value = "  -0.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.15%  "

# Row 43
$ws.Range("E43").Value = "  -2.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.84%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$ws.Range("E46").Value = "  -2.68%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0952"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.39%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.448.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.74%  "


# Row 48 - was FTXToken, now NEARProtocol
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "

# Row 49 - was NEARProtocol, now FTXToken
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.31%  "

# Row 51 - was ARBITRUM, now Celestia
$ws.Range("B51").Value = "Celestia"
$ws.Range("C51").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.64%  "
